$d = $word.ActiveDocument

# Locate the "Dheeraj Chand" name heading paragraph.
$found = $d.Content.Find.Execute("Dheeraj Chand", $true, $false, $false, $false,
                                  $false, $true, 1, $false, "", 0)

$namePara = $d.Paragraphs(1)
$nameRange = $namePara.Range.Duplicate

# Pull this paragraph's own WordprocessingML so we can re-emit it unchanged
# alongside the new paragraph. (Range.InsertXML *replaces* the range it is
# called on, so to insert a sibling paragraph right after this one without
# disturbing the paragraph that currently follows it, we replace this
# paragraph's range with "itself + the new paragraph".)
$full = $nameRange.WordOpenXML
$bodyStart = $full.IndexOf("<w:body>") + 8
$firstParaEnd = $full.IndexOf("</w:p>", $bodyStart) + 6
$nameParaXml = $full.Substring($bodyStart, $firstParaEnd - $bodyStart)

# Word tags the re-emitted paragraph with fresh w14:paraId / rsid scaffolding
# attributes; strip those back off so the paragraph round-trips unchanged.
$nameParaXml = $nameParaXml -replace '^<w:p\b[^>]*>', '<w:p>'

# The new, centered, unformatted contact-info paragraph (no inherited
# bold / 28pt run formatting from the name line).
$contactText = "202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX"
$contactParaXml = "<w:p><w:pPr><w:jc w:val=""center""/></w:pPr><w:r><w:t>$contactText</w:t></w:r></w:p>"

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' + $nameParaXml + $contactParaXml + '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$nameRange.InsertXML($xml)
